$d = $word.ActiveDocument

# 1 & 2) Remove the hyperlinks wrapping "Non-standard evaluation" and
#        "critical for correct calculation", keeping their visible text
#        (this is what Hyperlink.Delete does: it removes the hyperlink
#        field/formatting but leaves the display text behind as a plain run).
$targets = @('"Non-standard evaluation"', 'critical for correct calculation')
foreach ($target in $targets) {
    for ($i = $d.Hyperlinks.Count; $i -ge 1; $i--) {
        $h = $d.Hyperlinks.Item($i)
        if ($h.TextToDisplay -eq $target) {
            $h.Delete()
        }
    }
}

# 3) "critical for correct calculation" was split across 3 runs
#    ("critical for ", "c", "orrect calculation") even after unlinking above.
#    Re-assert the text via Find & Replace so it collapses into a single run
#    while keeping the existing (now un-hyperlinked) character formatting.
$d.Content.Find.Execute("critical for correct calculation", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "critical for correct calculation", 2) | Out-Null

# 4) "block-if/else statement" was split across 3 runs
#    ("block-if/else state", "m", "ent"). Collapse into a single run the same way.
$d.Content.Find.Execute("block-if/else statement", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "block-if/else statement", 2) | Out-Null
